$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Random4")

# --- Table 1 (rows 3-8): fill in the ratio formulas (I..M) that were previously blank ---
$ws2.Range("I3:I6").Formula = "=E3/C3"
$ws2.Range("J3:J6").Formula = "=D3/C3"
$ws2.Range("K3:K6").Formula = "=F3/C3"
$ws2.Range("L3:L6").Formula = "=I3/K3"
$ws2.Range("M3:M6").Formula = "=J3/K3"

$ws2.Range("I7:I8").Formula = "=E7/C7"
$ws2.Range("J7:J8").Formula = "=D7/C7"
$ws2.Range("K7:K8").Formula = "=F7/C7"
$ws2.Range("L7:L8").Formula = "=I7/K7"
$ws2.Range("M7:M8").Formula = "=J7/K7"

# --- Table 2 (rows 15-20): new benchmark run results for depths 1-4 (rows 15-18) ---
$ws2.Range("D15").Value = 103529
$ws2.Range("E15").Value = 1103921
$ws2.Range("F15").Value = 167078359
$ws2.Range("G15").Value = 15094
$ws2.Range("H15").Value = 76776

$ws2.Range("D16").Value = 216248
$ws2.Range("E16").Value = 2535625
$ws2.Range("F16").Value = 81413133
$ws2.Range("G16").Value = 19173
$ws2.Range("H16").Value = 76940

$ws2.Range("D17").Value = 369062
$ws2.Range("E17").Value = 4370251
$ws2.Range("F17").Value = 60425721
$ws2.Range("G17").Value = 30969
$ws2.Range("H17").Value = 126820

$ws2.Range("D18").Value = 1587241
$ws2.Range("E18").Value = 18981224
$ws2.Range("F18").Value = 19499363
$ws2.Range("G18").Value = 33687
$ws2.Range("H18").Value = 140860

# Depths 5-6 (rows 19-20) reuse the prior depth-5/6 results (no new run yet) - clear the ratio formulas
$ws2.Range("D19").Value = 1144277
$ws2.Range("E19").Value = 13574233
$ws2.Range("F19").Value = 9309472
$ws2.Range("G19").Value = 51209
$ws2.Range("H19").Value = 159256
$ws2.Range("I19:M20").ClearContents()

$ws2.Range("D20").Value = 3271445
$ws2.Range("E20").Value = 38366544
$ws2.Range("F20").Value = 1942471
$ws2.Range("G20").Value = 53842
$ws2.Range("H20").Value = 153196

# Code version labels: rows 15-18 are the new 2.3.13 run, rows 19-20 still tagged as 2.3.12
$ws2.Range("AB15:AB18").Value = "2.3.13"
$ws2.Range("AB19:AB20").Value = "2.3.12"

# --- New summary row 22 mirroring row 10's totals, for the updated table ---
$ws2.Range("G22").Formula = "=SUM(G15:G20)"
$ws2.Range("H22").Formula = "=SUM(H15:H20)"
$ws2.Range("M22").Formula = "=SUM(M15:M20)"
$ws2.Range("N22").Formula = "=SUM(N15:N20)/6"
$ws2.Range("O22").Formula = "=SUM(O15:O20)/6"
$ws2.Range("P22:Y22").Formula = "=SUM(P15:P20)/6"

# --- Selections (view state) ---
$ws1 = $wb.Worksheets.Item("Random7")
$ws1.Range("I8:M8").Select()
$ws2.Range("I7:M8").Select()
